# Portfolio: Updated Items data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (DONNY OSMOND) - "development" wording tweak + "highlights" year added
$ws.Range("I2").Value = "Layout aesthetic; content building; creating the navigation of each category and collection as well as supporting items; callouts; content development; detail to user navigation from start to end"
$ws.Range("J2").Value = "Donny Osmond himself expressed his gratitude to me personally for the vision created that represents his brand at the 2018 Las Vegas Market Show!"

# Row 3 (SCOTT LIVING) - "development" wording tweak
$ws.Range("I3").Value = "Layout aesthetic + developing working files, complying with the developing entity of SL / SBG. Collecting data, adjusting studio photography, page and spread layouts, cover + developing signature data for print."

# Row 4 (CF PRINT MEDIA) - "development" copy edits (artist->artists, punctuation, blank lines collapsed, br tag reformatted)
$ws.Range("I4").Value = "<p>TOP: create a fun and engaging visual informational campaign that included a brochure, stand, poster, and email blast for PROP 65. Explored visual engagement to take the viewer on a tour of the complex topic.</p>`n<p>BOTTOM: Worked with freelance artists to create and manipulate characters. Developed an array of patterns in Illustrator displayed on “giveaways” such as recycled bags, umbrellas, and other media.</p>`n<p>`n    Created mocks and hangtag.<br />`n    Samples available.`n</p>"

# Row 5 (CF WEB MEDIA) - "showcased" gets trailing period, "development" copy edits (end-user hyphenation, rewrites)
$ws.Range("F5").Value = "Oops! Correction to an email campaign sent with broken links. Also shown, web slides, and digital web assets."
$ws.Range("I5").Value = "<p>TOP: communicate with the audience about an error from the original email. The goal was to present something fun and convey accurate information to the end-user minimizing frustration and acknowledging the mistake. </p>`n<p>2: Web slides hosted on one of the Coaster Company websites. Based on the time of year, I created slides that told a story of a particular event: market show announcements, social media links, special promotions. </p>`n<p>3: Develop a visual web portal that allows end-user to get access to web-friendly converted print material: catalogs, images, and informational spreadsheets.</p>"

# Row 6 (PERENNIAL STRENGTH) - "development" rewritten as HTML paragraphs
$ws.Range("I6").Value = "<p>Create an identity for a Cross-Fit coach. Developed color scheme, logo, social media banners, mailers, business cards, and website visual aspects and navigation. Worked with a web developer for live website:</p>`n<p>visit: perennialstrength.com</p>"

# The multi-line replacement text can trigger row auto-fit; keep the original fixed row heights
$ws.Rows.Item(2).RowHeight = 18
$ws.Rows.Item(3).RowHeight = 18
$ws.Rows.Item(4).RowHeight = 18
$ws.Rows.Item(5).RowHeight = 18
$ws.Rows.Item(6).RowHeight = 18

# Restore the cell selection that was active when the author saved the file
$ws.Range("J8").Select()
